# Add five more "Multiple Choice" answer paragraphs (B, C, D, A, C) after
# the existing last answer, matching the ListParagraph style / numbering
# used by the rest of the answer key, and move the "_GoBack" bookmark so
# it ends up right after the text of the new final paragraph (mirroring
# where Word leaves it after the last edit).

$d = $word.ActiveDocument

$newAnswers = @("B", "C", "D", "A", "C")
foreach ($answer in $newAnswers) {
    $lastPara = $d.Paragraphs.Last
    $r = $lastPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.InsertAfter($answer)
}

# The "_GoBack" bookmark currently sits at the end of the old last
# paragraph; remove it so it can be recreated at the new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Temporarily append a placeholder character to the final paragraph. This
# sidesteps a quirk where adding a bookmark collapsed exactly on the
# paragraph-mark position snaps to the wrong spot; with the placeholder in
# place, the target offset is a normal text position instead.
$finalPara = $d.Paragraphs.Last
$fr = $finalPara.Range
$fr.Collapse(0)
$fr.InsertAfter("X")

$finalPara = $d.Paragraphs.Last
$markerPos = $finalPara.Range.End - 2
$bookmarkRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the placeholder character again; the bookmark stays put.
$placeholderRange = $d.Range($finalPara.Range.End - 2, $finalPara.Range.End - 1)
$placeholderRange.Delete()
